# aggiornamento fino a 27/05
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44330, 1, 1, 16.63616702711695),
    @(44331, 0, 1, 16.63616702711695),
    @(44332, 0, 1, 16.63616702711695),
    @(44333, 0, 1, 16.63616702711695),
    @(44334, 3, 4, 66.54466810846782),
    @(44335, 0, 4, 66.54466810846782),
    @(44336, 0, 4, 66.54466810846782),
    @(44337, 0, 3, 49.90850108135086),
    @(44338, 0, 3, 49.90850108135086),
    @(44339, 0, 3, 49.90850108135086),
    @(44340, 0, 3, 49.90850108135086),
    @(44341, 0, 0, 0),
    @(44342, 0, 0, 0),
    @(44343, 0, 0, 0)
)

$startRow = 256
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

# copy the style from the last existing data row (255) onto the newly added rows
$srcRow = $startRow - 1
$endRow = $startRow + $data.Count - 1
$ws.Range("A$srcRow`:D$srcRow").Copy() | Out-Null
$ws.Range("A$startRow`:D$endRow").PasteSpecial(-4122) | Out-Null
